$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new reference row (A21) with the new URL, matching the style
# used by the other hyperlink cells, and register the hyperlink itself.
$newUrl = "https://codes-sources.commentcamarche.net/source/53687-lister-fichiers-et-repertoires-multiplateforme"

$ws.Hyperlinks.Add($ws.Range("A21"), $newUrl)
$ws.Range("A21").Style = "Lien hypertexte"

# Move the active selection to A22, just past the newly added row.
$ws.Range("A22").Select()
